# Auto-generated script applying the values from the commit diff
# (scheduled runner updating cached market-board price snapshots)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2003.8
$ws.Range("J40").Value = 2219
$ws.Range("L40").Value = 2219
$ws.Range("N40").Value = -2569
$ws.Range("H88").Value = 10579.667
$ws.Range("J88").Value = 10579.667
$ws.Range("L88").Value = 10579.667
$ws.Range("N88").Value = -11391.667
$ws.Range("H91").Value = 10579.667
$ws.Range("J91").Value = 10579.667
$ws.Range("L91").Value = 10579.667
$ws.Range("N91").Value = -13387.667
$ws.Range("H112").Value = 16509.5
$ws.Range("J112").Value = 26021.2
$ws.Range("L112").Value = 78063.60000000001
$ws.Range("N112").Value = -80279.60000000001
$ws.Range("H129").Value = 2235.6667
$ws.Range("I129").Value = 2235.6667
$ws.Range("K129").Value = 6707.000100000001
$ws.Range("M129").Value = -1707.000100000001
$ws.Range("H132").Value = 16351.404
$ws.Range("I132").Value = 16351.404
$ws.Range("K132").Value = 49054.212
$ws.Range("M132").Value = -46524.212

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 341160.16
$ws.Range("I74").Value = 547023.25
$ws.Range("J74").Value = 17661
$ws.Range("K74").Value = 547023.25
$ws.Range("L74").Value = 17661
$ws.Range("M74").Value = -546149.25
$ws.Range("N74").Value = -19409
$ws.Range("H77").Value = 341160.16
$ws.Range("I77").Value = 547023.25
$ws.Range("J77").Value = 17661
$ws.Range("K77").Value = 2735116.25
$ws.Range("L77").Value = 88305
$ws.Range("M77").Value = -2730748.25
$ws.Range("N77").Value = -97041
$ws.Range("H88").Value = 7291.4443
$ws.Range("J88").Value = 11559.4
$ws.Range("L88").Value = 11559.4
$ws.Range("N88").Value = -12371.4
$ws.Range("H91").Value = 7291.4443
$ws.Range("J91").Value = 11559.4
$ws.Range("L91").Value = 11559.4
$ws.Range("N91").Value = -14367.4
$ws.Range("H102").Value = 2616.4666
$ws.Range("I102").Value = 2488.2856
$ws.Range("J102").Value = 4411
$ws.Range("K102").Value = 2488.2856
$ws.Range("L102").Value = 4411
$ws.Range("M102").Value = -866.2856000000002
$ws.Range("N102").Value = -7655

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 7694267.5
$ws.Range("J94").Value = 25002574
$ws.Range("L94").Value = 25002574
$ws.Range("N94").Value = -25003476
$ws.Range("H105").Value = 1672.6086
$ws.Range("J105").Value = 2218.4
$ws.Range("L105").Value = 2218.4
$ws.Range("N105").Value = -5712.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7143710.5
$ws.Range("I31").Value = 7143710.5
$ws.Range("K31").Value = 7143710.5
$ws.Range("M31").Value = -7143415.5
$ws.Range("H34").Value = 7143710.5
$ws.Range("I34").Value = 7143710.5
$ws.Range("K34").Value = 7143710.5
$ws.Range("M34").Value = -7143508.5
$ws.Range("H62").Value = 5076.273
$ws.Range("J62").Value = 4957.5386
$ws.Range("L62").Value = 4957.5386
$ws.Range("N62").Value = -6205.5386
$ws.Range("H65").Value = 5076.273
$ws.Range("J65").Value = 4957.5386
$ws.Range("L65").Value = 24787.693
$ws.Range("N65").Value = -31027.693
$ws.Range("H103").Value = 13899.6
$ws.Range("I103").Value = 13899.6
$ws.Range("K103").Value = 13899.6
$ws.Range("M103").Value = -12727.6
$ws.Range("H132").Value = 60173.06
$ws.Range("I132").Value = 63621.375
$ws.Range("K132").Value = 190864.125
$ws.Range("M132").Value = -188334.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 183.35715
$ws.Range("I2").Value = 181.18182
$ws.Range("K2").Value = 1087.09092
$ws.Range("M2").Value = -974.0909199999999
$ws.Range("H34").Value = 5166.722
$ws.Range("I34").Value = 2001
$ws.Range("K34").Value = 6003
$ws.Range("M34").Value = -5919
$ws.Range("H39").Value = 17855.715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 12330
$ws.Range("I41").Value = 2260
$ws.Range("J41").Value = 22400
$ws.Range("K41").Value = 2260
$ws.Range("L41").Value = 22400
$ws.Range("M41").Value = -1905
$ws.Range("N41").Value = -23110
$ws.Range("H80").Value = 5285.037
$ws.Range("J80").Value = 8334.111000000001
$ws.Range("L80").Value = 8334.111000000001
$ws.Range("N80").Value = -10330.111
$ws.Range("H83").Value = 5285.037
$ws.Range("J83").Value = 8334.111000000001
$ws.Range("L83").Value = 41670.55500000001
$ws.Range("N83").Value = -51654.55500000001
$ws.Range("H113").Value = 2649.6667
$ws.Range("I113").Value = 1420.3334
$ws.Range("K113").Value = 1420.3334
$ws.Range("M113").Value = 749.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1344.238
$ws.Range("I22").Value = 1159.9231
$ws.Range("J22").Value = 1643.75
$ws.Range("K22").Value = 1159.9231
$ws.Range("L22").Value = 1643.75
$ws.Range("M22").Value = -864.9231
$ws.Range("N22").Value = -2233.75
$ws.Range("H27").Value = 1344.238
$ws.Range("I27").Value = 1159.9231
$ws.Range("J27").Value = 1643.75
$ws.Range("K27").Value = 1159.9231
$ws.Range("L27").Value = 1643.75
$ws.Range("M27").Value = -1052.9231
$ws.Range("N27").Value = -1857.75
$ws.Range("H42").Value = 14796.444
$ws.Range("I42").Value = 30099
$ws.Range("J42").Value = 12883.625
$ws.Range("K42").Value = 30099
$ws.Range("L42").Value = 12883.625
$ws.Range("M42").Value = -29536
$ws.Range("N42").Value = -14009.625
$ws.Range("H49").Value = 14796.444
$ws.Range("I49").Value = 30099
$ws.Range("J49").Value = 12883.625
$ws.Range("K49").Value = 30099
$ws.Range("L49").Value = 12883.625
$ws.Range("M49").Value = -29952
$ws.Range("N49").Value = -13177.625
$ws.Range("H82").Value = 2099
$ws.Range("I82").Value = 1987.7778
$ws.Range("J82").Value = 2432.6667
$ws.Range("K82").Value = 1987.7778
$ws.Range("L82").Value = 2432.6667
$ws.Range("M82").Value = -1626.7778
$ws.Range("N82").Value = -3154.6667
$ws.Range("H85").Value = 2099
$ws.Range("I85").Value = 1987.7778
$ws.Range("J85").Value = 2432.6667
$ws.Range("K85").Value = 1987.7778
$ws.Range("L85").Value = 2432.6667
$ws.Range("M85").Value = -739.7778000000001
$ws.Range("N85").Value = -4928.6667
$ws.Range("H93").Value = 1589.875
$ws.Range("I93").Value = 1194.7273
$ws.Range("J93").Value = 2459.2
$ws.Range("K93").Value = 1194.7273
$ws.Range("L93").Value = 2459.2
$ws.Range("M93").Value = 53.27269999999999
$ws.Range("N93").Value = -4955.2
$ws.Range("H136").Value = 3169.4062
$ws.Range("I136").Value = 2560.3333
$ws.Range("K136").Value = 7680.999899999999
$ws.Range("M136").Value = -5130.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 20749.5
$ws.Range("I81").Value = 20749.5
$ws.Range("K81").Value = 41499
$ws.Range("M81").Value = -40438
$ws.Range("H84").Value = 20749.5
$ws.Range("I84").Value = 20749.5
$ws.Range("K84").Value = 207495
$ws.Range("M84").Value = -202191
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 954.94446
$ws.Range("I107").Value = 955.4545000000001
$ws.Range("K107").Value = 2866.3635
$ws.Range("M107").Value = -946.3635000000004
$ws.Range("H122").Value = 52246.414
$ws.Range("I122").Value = 57794.27
$ws.Range("J122").Value = 4165
$ws.Range("K122").Value = 173382.81
$ws.Range("L122").Value = 12495
$ws.Range("M122").Value = -170932.81
$ws.Range("N122").Value = -17395
$ws.Range("H136").Value = 21649.275
$ws.Range("I136").Value = 28396.762
$ws.Range("J136").Value = 3937.125
$ws.Range("K136").Value = 85190.28599999999
$ws.Range("L136").Value = 11811.375
$ws.Range("M136").Value = -82640.28599999999
$ws.Range("N136").Value = -16911.375
